$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The authoritative commit adds one new daily price observation for Perejil
# (Vega Central Mapocho de Santiago) dated serial 44447, inserted above the
# existing row 108. This pushes every subsequent data row (108-203) down by
# one (new dimension A1:R204). Inserting a whole row keeps all the existing
# data/styles intact and automatically extends the used range/dimension.
$ws.Rows("108:108").Insert()

$ws.Range("A108").Value = 9
$ws.Range("B108").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C108").Value = "Metropolitana"
$ws.Range("D108").Value = 44447
$ws.Range("E108").Value = 13
$ws.Range("F108").Value = 100112044
$ws.Range("G108").Value = "Perejil"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 97
$ws.Range("K108").Value = 10000
$ws.Range("L108").Value = 11000
$ws.Range("M108").Value = 10505
$ws.Range("N108").Value = "`$/docena de atados"
$ws.Range("O108").Value = "Región Metropolitana"
$ws.Range("P108").Value = 3502
$ws.Range("Q108").Value = 3
$ws.Range("R108").Value = "Hortaliza"
